# Atividades.xlsx update: refresh the activity rows with updated codes/titles
# and add a third activity row (commit: "Fixed some bugs in auditing").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: Atividade01 ---
$ws.Range("A1").Value = 158070
$ws.Range("A1").Style = "Normal"
$ws.Range("B1").Value = "Atividade01_atualizado"
$ws.Range("B1").WrapText = $true
$ws.Range("C1").ClearContents()

# --- Row 2: Atividade02 ---
$ws.Range("A2").Value = 158071
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Atividade02_atualizado"
$ws.Range("C2").ClearContents()

# --- Row 3 (new): Atividade03 ---
$ws.Range("A3").Value = 158072
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Atividade03_atualizado"

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 13.8

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 68.2
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666

# --- Selection matches the post-edit cursor position ---
$ws.Range("B3").Select() | Out-Null
